$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The anonymization rewrite dropped one duplicate store row (old row 29, "Лента")
# which shifts rows 30-31 up to 29-30, shrinking the used range to A1:B30.
$ws.Rows.Item(29).Delete()

# Replace the lat/lon coordinate pairs in column B with the anonymized district names.
$ws.Range("B2").Value = 'Колпинский район'
$ws.Range("B3").Value = 'Московский район'
$ws.Range("B4").Value = 'Центральный район'
$ws.Range("B5").Value = 'Всеволожский район (Ленинградская область)'
$ws.Range("B6").Value = 'Центральный район'
$ws.Range("B7").Value = 'Центральный район'
$ws.Range("B8").Value = 'Выборгский район'
$ws.Range("B9").Value = 'Фрунзенский район'
$ws.Range("B10").Value = 'Выборгский район'
$ws.Range("B11").Value = 'Красногвардейский район'
$ws.Range("B12").Value = 'Красногвардейский район'
$ws.Range("B13").Value = 'Всеволожский район (Ленинградская область)'
$ws.Range("B14").Value = 'Приморский район'
$ws.Range("B15").Value = 'Фрунзенский район'
$ws.Range("B16").Value = 'Центральный район'
$ws.Range("B17").Value = 'Красногвардейский район'
$ws.Range("B18").Value = 'Красногвардейский район'
$ws.Range("B19").Value = 'Фрунзенский район'
$ws.Range("B20").Value = 'Центральный район'
$ws.Range("B21").Value = 'Красногвардейский район'
$ws.Range("B22").Value = 'Невский район'
$ws.Range("B23").Value = 'Центральный район'
$ws.Range("B24").Value = 'Выборгский район'
$ws.Range("B25").Value = 'Всеволожский район (Ленинградская область)'
$ws.Range("B26").Value = 'Адмиралтейский район'
$ws.Range("B27").Value = 'Красногвардейский район'
$ws.Range("B28").Value = 'Приморский район'
$ws.Range("B29").Value = 'Московский район'
$ws.Range("B30").Value = 'Центральный район'
